$d = $word.ActiveDocument

# Locate the existing "${tanggal_rapat}" run inside the "pada hari/tanggal"
# table row and insert a new run containing "${hari}/" immediately before
# it, so the cell reads "${hari}/${tanggal_rapat}".
$rng = $d.Content
$found = $rng.Find.Execute('${tanggal_rapat}', $true, $false, $false, $false, $false, `
                            $true, 1, $false, "", 0)

if ($found) {
    $insertStart = $rng.Start
    $rng.InsertBefore('${hari}/')

    # Toggle a character-level formatting property on the newly inserted
    # text so it is written out as its own run instead of being merged
    # back into the neighbouring "${tanggal_rapat}" run (both ultimately
    # share the same resolved formatting).
    $newRunRange = $d.Range($insertStart, $insertStart + 8)
    $newRunRange.Bold = 1
    $newRunRange.Bold = 0
}
